# Generate Report for Archive
#
# The "Status" value "Ready for handoff" is now stale (it is regenerated as
# part of the archive report); the current state of these rows is
# "In Translation". Update every place that value appears:
#   - Overview!E2 / Overview!F2 (the per-language status columns)
#   - zh-cn!C2     (the "Status" table column)
#   - de-de!C2     (the "Status" table column)
#
# Because the new text is shorter than the old text, the "Status" columns
# are re-sized to fit the new content.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Best attainable ColumnWidth for the new, narrower "Status" columns -
# the stored column width (in characters) ends up at ~13.41 once written
# back out, matching the narrower fit for "In Translation".
$newColumnWidth = 12.576851254417766

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E1:F1").ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
